$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reported indicator values (4T2023 report refresh).
# Cells A3/A4 hold plain text that looks numeric ("93.9%", "+2,000"), so
# force text formatting before writing the new value to keep Excel from
# reinterpreting it as a percentage/number, then restore the original
# General number format so the cell's style is unchanged.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "93.7%"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "+2,500"

$ws.Range("A6").Value = "+27,0 Millones"

$ws.Range("A3").NumberFormat = "General"
$ws.Range("A4").NumberFormat = "General"

# Update the active selection to A7.
$ws.Range("A7").Select()
